# Fill in the newly-reported 2025 Q3 figures (row 44) that were added to
# the IK Konjunktur + Destatis/HWWI dashboard data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C44").Value2 = 198.46
$ws.Range("D44").Value2 = 155.97
$ws.Range("E44").Value2 = 154.01
$ws.Range("F44").Value2 = 345.47

$ws.Range("O44").Value2 = 509.66666666666669
$ws.Range("P44").Value2 = 90814.666666666686
$ws.Range("Q44").Value2 = 34594000
$ws.Range("R44").Value2 = 1121252666
$ws.Range("S44").Value2 = 6339248333
$ws.Range("T44").Value2 = 3419494000
$ws.Range("U44").Value2 = 2919754334
$ws.Range("V44").Value2 = 1683136667
$ws.Range("W44").Value2 = 1236617667

# The author's new selection after entering this row's data landed on the
# next (still partially empty) row.
$ws.Range("O45:W45").Select() | Out-Null

# Cosmetic re-save artifacts that a real Excel session also produces when
# the workbook is opened and saved again (reverts to Excel's default page
# margins and scrolls the view back to the top-left).
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
